# Apply BOM updates to the "ESC" worksheet:
#  - Insert a new component row (STLINK-V3SET) above the old ST-LINK/V2 row
#  - Zero out the ST-LINK/V2 quantity (replaced by the STLINK-V3SET)
#  - Insert 4 more rows in the blank area and populate 6 new components
#  - Re-select the "ESC" sheet as the active tab

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ESC")

# --- Insert the STLINK-V3SET row above the existing ST-LINK/V2 row (row 8) ---
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = "STLINK-V3SET"
$ws.Range("B8").Value = "497-18216-ND"
$ws.Range("C8").Value = 35
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 35
$ws.Range("F8").Value = "DEBUGGER/PROGRAMMER STM8 STM32"
$ws.Range("F8").HorizontalAlignment = -4131
$ws.Range("F8").VerticalAlignment = -4160
$ws.Range("F8").WrapText = $true
$ws.Range("G8").Value = "https://www.digikey.com/product-detail/en/stmicroelectronics/STLINK-V3SET/497-18216-ND/9636028"

# The old ST-LINK/V2 programmer (now row 9) is no longer used - qty -> 0
$ws.Range("D9").Value = 0

# --- Make room for 6 more new parts in the blank area (old rows 12-18, now 13-19) ---
$ws.Range("A13:A16").EntireRow.Insert()

# MIC5319-3.3YD5-TR  (LDO 5V->3.3V)
$ws.Range("A13").Value = "MIC5319-3.3YD5-TR "
$ws.Range("B13").Value = "576-1910-1-ND"
$ws.Range("C13").Value = 1.32
$ws.Range("D13").Value = 5
$ws.Range("E13").Formula = "=C13*D13"
$ws.Range("F13").Value = "LDO for 5V to 3.3V"
$ws.Range("G13").Value = "https://www.digikey.com/product-detail/en/microchip-technology/MIC5319-3-3YD5-TR/576-1910-1-ND/1799517"

# LD1117S50TR (12V->5V regulator)
$ws.Range("A14").Value = "LD1117S50TR"
$ws.Range("B14").Value = "497-6447-1-ND"
$ws.Range("C14").Value = 0.351
$ws.Range("D14").Value = 10
$ws.Range("E14").Formula = "=C14*D14"
$ws.Range("F14").Value = "12V to 5V regulator"
$ws.Range("G14").Value = "https://www.digikey.com/product-detail/en/stmicroelectronics/LD1117S50TR/497-6447-1-ND/1865475"

# FTSH-105-01-L-DV-K-TR (ARM10 connector)
$ws.Range("A15").Value = "FTSH-105-01-L-DV-K-TR"
$ws.Range("B15").Value = "SAM13165CT-ND"
$ws.Range("C15").Value = 3.6
$ws.Range("D15").Value = 2
$ws.Range("E15").Formula = "=C15*D15"
$ws.Range("F15").Value = "ARM10 connector"
$ws.Range("G15").Value = "https://www.digikey.com/products/en?keywords=SAM13165CT-ND"

# B3F-1020 (tactile switch)
$ws.Range("A16").Value = "B3F-1020"
$ws.Range("B16").Value = "SW402-ND"
$ws.Range("C16").Value = 0.32
$ws.Range("D16").Value = 5
$ws.Range("E16").Formula = "=C16*D16"
$ws.Range("F16").Value = "SWITCH TACTILE SPST-NO 0.05A 24V"
$ws.Range("G16").Value = "https://www.digikey.com/product-detail/en/omron-electronics-inc-emc-div/B3F-1020/SW402-ND/44059"

# STPS2L30A (schottky diode)
$ws.Range("A17").Value = "STPS2L30A"
$ws.Range("B17").Value = "497-3759-1-ND"
$ws.Range("C17").Value = 0.32
$ws.Range("D17").Value = 10
$ws.Range("E17").Formula = "=C17*D17"
$ws.Range("F17").Value = "Diode Schottky 30V 2A Surface Mount SMA"
$ws.Range("G17").Value = "https://www.digikey.com/product-detail/en/stmicroelectronics/STPS2L30A/497-3759-1-ND/691929"

# EEE-1VA100WR (10uF 35V electrolytic cap)
$ws.Range("A18").Value = "EEE-1VA100WR"
$ws.Range("B18").Value = "PCE3948CT-ND"
$ws.Range("C18").Value = 0.33
$ws.Range("D18").Value = 5
$ws.Range("E18").Formula = "=C18*D18"
$ws.Range("F18").Value = "10µF 35V Aluminum Electrolytic Capacitors Radial, Can - SMD  1000 Hrs @ 85°C"
$ws.Range("G18").Value = "https://www.digikey.com/product-detail/en/panasonic-electronic-components/EEE-1VA100WR/PCE3948CT-ND/766324"

# Put the focus/selection back on the ESC sheet near the new data
$ws.Activate()
$ws.Range("G22").Select()
